# Données_groupe_01 — data correction pass
# (commit: "modified data (because there were some logic problems)")
#
# The source data in columns A (ratio) and C (total) had logic problems
# for rows 13-61 (excluding the handful of rows that were already correct)
# and were recomputed / corrected. Column B is untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("données01")

# Best-effort: reflect the enlarged Excel window from the author's session.
# (Cosmetic view state; harmless no-op if the host doesn't persist it.)
$win = $excel.ActiveWindow
$win.Width = 25800
$win.Height = 13200

# Corrected values: cell address -> new value
$values = [ordered]@{
    "A13" = 32.32
    "C13" = 68
    "A14" = 8.91
    "C14" = 85
    "A16" = 16.72
    "C16" = 50
    "A17" = 23.59
    "C17" = 80
    "A18" = 22.759999999999998
    "C18" = 42
    "A19" = 5.96
    "C19" = 44
    "A20" = 28.23
    "C20" = 41
    "A21" = 8.129999999999999
    "C21" = 82
    "A22" = 23.9
    "C22" = 85
    "A23" = 6.6199999999999992
    "C23" = 84
    "A24" = 43.87
    "C24" = 79
    "A25" = 33.31
    "C25" = 50
    "A26" = 41.47
    "C26" = 75
    "A27" = 9.19
    "C27" = 81
    "A28" = 9.0399999999999991
    "C28" = 79
    "A29" = 3.17
    "C29" = 87
    "A30" = 38.279999999999994
    "C30" = 74
    "A31" = 83.71
    "C31" = 91
    "A32" = 34
    "C32" = 80
    "A33" = 3.1
    "C33" = 36
    "A34" = 26.35
    "C34" = 86
    "A35" = 4.1099999999999994
    "C35" = 91
    "A36" = 20.91
    "C36" = 89
    "A38" = 6.1400000000000006
    "C38" = 91
    "A39" = 20.599999999999998
    "C39" = 77
    "A40" = 30.270000000000003
    "C40" = 66
    "A42" = 3.2
    "C42" = 90
    "A43" = 29.659999999999997
    "C43" = 91
    "A44" = 4.3600000000000003
    "C44" = 73
    "A45" = 7.9
    "C45" = 84
    "A47" = 23.880000000000003
    "C47" = 90
    "A48" = 2.48
    "C48" = 81
    "A50" = 17.34
    "C50" = 83
    "A51" = 19.53
    "C51" = 80
    "A52" = 9.5299999999999994
    "C52" = 66
    "A54" = 11.51
    "C54" = 77
    "A55" = 4.9799999999999995
    "C55" = 89
    "A56" = 23.98
    "C56" = 67
    "A57" = 46.239999999999995
    "C57" = 75
    "A58" = 44.26
    "C58" = 77
    "A59" = 10.620000000000001
    "C59" = 49
    "A60" = 24.279999999999998
    "C60" = 62
    "A61" = 11.83
    "C61" = 82
}

foreach ($cell in $values.Keys) {
    $ws.Range($cell).Value2 = $values[$cell]
}
